$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175147294998169
$ws.Range("B1").Value = 2.15348744392395
$ws.Range("C1").Value = 3.577345371246338
$ws.Range("D1").Value = 3.44361686706543
$ws.Range("E1").Value = 1.16057288646698
